# Update countries & provincias Spain
# - Reorders Niger / Sri Lanka so that Sri Lanka's (updated) row comes first,
#   followed by Niger's previous row of data (shifted down one row).
# - Reorders El Salvador / San Marino / Mali so El Salvador's (updated) row
#   comes first, followed by San Marino's and Mali's previous rows of data
#   (each shifted down one row).
# - Refreshes the "last updated" timestamp string.
# - Applies new case totals for several countries (USA, Germany, Mexico,
#   Pakistan, Australia, Sri Lanka, Niger, El Salvador, San Marino, Mali,
#   Mongolia).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country name swaps (rows keep their statistics row position, but the
#     country label attached to that row changes because the country that
#     has fresh data moves to the earlier row) ---

# Row 99 was Niger -> becomes Sri Lanka (fresh data)
# Row 100 was Sri Lanka -> becomes Niger (old Niger numbers, shifted down)
$ws.Range("A99").Value = "Sri Lanka"
$ws.Range("A100").Value = "Niger"

# Row 112 was San Marino -> becomes El Salvador (fresh data)
# Row 113 was Mali -> becomes San Marino (old San Marino numbers, shifted down)
# Row 114 was El Salvador -> becomes Mali (old Mali numbers, shifted down)
$ws.Range("A112").Value = "El Salvador"
$ws.Range("A113").Value = "San Marino"
$ws.Range("A114").Value = "Mali"

# --- Updated "last refreshed" footer text ---
$ws.Range("A1").Value = "Datos actualizados a 5 de Mayo de 2020 a las 06:03"

# --- Numeric data updates ---

# Row 4: Estados Unidos
$ws.Range("B4").Value = 1212900
$ws.Range("C4").Value = 65
$ws.Range("D4").Value = 188068
$ws.Range("E4").Value = 954911

# Row 9: Alemania
$ws.Range("D9").Value = 135100
$ws.Range("E9").Value = 24059

# Row 24: Mexico
$ws.Range("D24").Value = 15938
$ws.Range("E24").Value = 6696

# Row 27: Pakistan
$ws.Range("B27").Value = 21501
$ws.Range("C27").Value = 560
$ws.Range("D27").Value = 5782
$ws.Range("E27").Value = 15233
$ws.Range("G27").Value = 10
$ws.Range("H27").Value = 486

# Row 51: Australia
$ws.Range("E51").Value = 865
$ws.Range("G51").Value = 1
$ws.Range("H51").Value = 96

# Row 99: now Sri Lanka (fresh numbers)
$ws.Range("B99").Value = 755
$ws.Range("C99").Value = 4
$ws.Range("D99").Value = 194
$ws.Range("E99").Value = 553
$ws.Range("F99").Value = 1
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 8

# Row 100: now Niger (old numbers shifted down from the previous row 99)
$ws.Range("B100").Value = 755
$ws.Range("C100").Value = 0
$ws.Range("D100").Value = 534
$ws.Range("E100").Value = 184
$ws.Range("F100").Value = 0
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 37

# Row 112: now El Salvador (fresh numbers)
$ws.Range("B112").Value = 587
$ws.Range("C112").Value = 32
$ws.Range("D112").Value = 201
$ws.Range("E112").Value = 373
$ws.Range("F112").Value = 3
$ws.Range("G112").Value = 0
$ws.Range("H112").Value = 13

# Row 113: now San Marino (old numbers shifted down from the previous row 112)
$ws.Range("B113").Value = 582
$ws.Range("C113").Value = 0
$ws.Range("D113").Value = 86
$ws.Range("E113").Value = 455
$ws.Range("F113").Value = 5
$ws.Range("G113").Value = 0
$ws.Range("H113").Value = 41

# Row 114: now Mali (old numbers shifted down from the previous row 113)
$ws.Range("B114").Value = 580
$ws.Range("C114").Value = 0
$ws.Range("D114").Value = 223
$ws.Range("E114").Value = 328
$ws.Range("F114").Value = 0
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = 29

# Row 175: Mongolia
$ws.Range("B175").Value = 41
$ws.Range("C175").Value = 1
$ws.Range("D175").Value = 29
$ws.Range("E175").Value = 29
